$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 21-23: column C renumbers (same text, but we just set values directly)
$ws.Range("C21").Value = "test7_1.csv"
$ws.Range("C22").Value = "test7_2.csv"
$ws.Range("C23").Value = "test7_3.csv"

# Row 24 (8.1): add C24
$ws.Range("C24").Value = "test7_1.csv"

# Row 25 (8.2): add C25
$ws.Range("C25").Value = "test7_2.csv"

# Row 26 (8.3): B26 text changes, add C26
$ws.Range("B26").Value = "VaR from Simulation -- compare to 8.2 values"
$ws.Range("C26").Value = "test7_2.csv"

# Row 27 (8.4): B27 stays same text, add C27
$ws.Range("B27").Value = "ES From Normal Distribution"
$ws.Range("C27").Value = "test7_1.csv"

# Row 28 (8.5): B28 stays same text, add C28
$ws.Range("B28").Value = "ES from T Distribution"
$ws.Range("C28").Value = "test7_2.csv"

# Row 29 (8.6): B29 text changes, add C29
$ws.Range("B29").Value = "ES from Simulation -- compare to 8.5 values"
$ws.Range("C29").Value = "test7_2.csv"

# Row 30 (9.1): B30 stays same text (renumbered only)
$ws.Range("B30").Value = "Gaussian Copula with a  Normal and T distribution"

# Row 31 (9.2): B31 stays same text (renumbered only)
$ws.Range("B31").Value = "Gaussian Copula with a  2 T distributions"

# Row 32 (10.1): B32 stays same text (renumbered only)
$ws.Range("B32").Value = "VaR/ES on 2 levels from simulated values"

# Update sheet view: scroll position + selection (B30, scrolled so row 10 is near top)
$ws.Activate()
$ws.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("B30").Select()
